$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = '28.519.63'
$ws.Range("E2").Value = '  -1.11%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = '1.870.20'
$ws.Range("E3").Value = '  -1.59%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -2.83%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "'315.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.87%  '

# Row 6: 'USDC' -> 'USDC'
$ws.Range("D6").Value = "'1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.38%  '

# Row 7: 'XRP' -> 'XRP'
$ws.Range("D7").Value = "'0.5088"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.23%  '

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").Value = "'0.3901"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.95%  '

# Row 9: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D9").Value = "'0.08352"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.51%  '

# Row 10: 'Polygon' -> 'Polygon'
$ws.Range("D10").Value = "'1.105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.10%  '

# Row 11: 'OKB' -> 'Polkadot'
$ws.Range("B11").Value = 'Polkadot'
$ws.Range("C11").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D11").Value = "'6.205"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.21%  '

# Row 12: 'Polkadot' -> 'WrappedEther'
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.864.47'
$ws.Range("E12").Value = '  +5.14%  '

# Row 13: 'WrappedEther' -> 'Solana'
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = "'20.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.16%  '

# Row 14: 'Solana' -> 'Chainlink'
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'7.254"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.03%  '

# Row 15: 'Chainlink' -> 'BinanceUSD'
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = "'1.009"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.99%  '

# Row 16: 'BinanceUSD' -> 'ShibaInu'
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.00001099"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.20%  '

# Row 17: 'ShibaInu' -> 'Litecoin'
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = "'91.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.14%  '

# Row 18: 'Litecoin' -> 'TRON'
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = "'0.06740"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.91%  '

# Row 19: 'TRON' -> 'Avalanche'
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = "'17.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.83%  '

# Row 20: 'Avalanche' -> 'Dai'
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = "'1.007"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.47%  '

# Row 21: 'Dai' -> 'Uniswap'
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'5.913"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.42%  '

# Row 22: 'Uniswap' -> 'WrappedBTC'
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '28.571.47'
$ws.Range("E22").Value = '  -0.96%  '

# Row 23: 'WrappedBTC' -> 'Cosmos'
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = "'11.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.89%  '

# Row 24: 'Cosmos' -> 'Toncoin'
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = "'2.206"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.88%  '

# Row 25: 'Toncoin' -> 'WrappedliquidstakedEther2.0'
$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.084.32'
$ws.Range("E25").Value = '  +4.84%  '

# Row 26: 'WrappedliquidstakedEther2.0' -> 'Monero'
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = "'156.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.20%  '

# Row 27: 'Monero' -> 'EthereumClassic'
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = "'20.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.39%  '

# Row 28: 'EthereumClassic' -> 'LidoDAOToken'
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = "'2.416"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.03%  '

# Row 29: 'LidoDAOToken' -> 'BitcoinCash'
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = "'126.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.44%  '

# Row 30: 'BitcoinCash' -> 'Stellar'
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = "'0.1037"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.86%  '

# Row 31: 'Stellar' -> 'ImmutableX'
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'1.042"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.56%  '

# Row 32: 'ImmutableX' -> 'Filecoin'
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'5.729"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.07%  '

# Row 33: 'Filecoin' -> 'HuobiToken'
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = "'3.620"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.76%  '

# Row 34: 'HuobiToken' -> 'VeChain'
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").Value = "'0.02453"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.64%  '

# Row 35: 'VeChain' -> 'Hedera'
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = "'0.06584"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.36%  '

# Row 36: 'Hedera' -> 'FraxShare'
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = "'8.921"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.28%  '

# Row 37: 'FraxShare' -> 'Algorand'
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = "'0.2159"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.61%  '

# Row 38: 'Algorand' -> 'InternetComputer(DFINITY)'
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = "'5.040"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.91%  '

# Row 39: 'InternetComputer(DFINITY)' -> 'ARBITRUM'
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = "'1.180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.27%  '

# Row 40: 'ARBITRUM' -> 'TrustWalletToken'
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'1.234"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.66%  '

# Row 41: 'TrustWalletToken' -> 'TheSandbox'
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'0.6354"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.90%  '

# Row 42: 'TheSandbox' -> 'Aptos'
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = "'11.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.77%  '

# Row 43: 'Aptos' -> 'Frax'
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = "'1.008"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.88%  '

# Row 44: 'Frax' -> 'Decentraland'
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = "'0.5991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.05%  '

# Row 45: 'Decentraland' -> 'EnergySwap'
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'13.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.19%  '

# Row 46: 'EnergySwap' -> 'PancakeSwap'
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = "'3.683"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.62%  '

# Row 47: 'PancakeSwap' -> 'NEARProtocol'
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'2.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.39%  '

# Row 48: 'NEARProtocol' -> 'EOS'
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").Value = "'1.214"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.24%  '

# Row 49: 'EOS' -> 'Quant'
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = "'122.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.53%  '

# Row 50: 'Quant' -> 'Cronos'
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.06804"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.74%  '

# Row 51: 'Cronos' -> 'Aave'
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = "'76.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.81%  '
